$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.165.30'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '2.376.13'
$ws.Range('E3').Value = '  +1.59%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''302.64'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('D6').Value = '''95.53'
$ws.Range('D7').Value = '''0.504'
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '''0.482'
$ws.Range('E9').Value = '  -2.60%  '
$ws.Range('D10').Value = '''34.37'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('E11').Value = '  +4.07%  '
$ws.Range('D12').Value = '''0.0788'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('D13').Value = '''18.26'
$ws.Range('E13').Value = '  -2.30%  '
$ws.Range('D14').Value = '''6.76'
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').Value = '2.740.26'
$ws.Range('E15').Value = '  +1.43%  '
$ws.Range('D16').Value = '2.393.82'
$ws.Range('E16').Value = '  +1.54%  '
$ws.Range('D17').Value = '''0.800'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').Value = '43.188.92'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').Value = '''11.97'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').Value = '0.0₃0889'
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('D22').Value = '''68.00'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').Value = '''235.69'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('E25').Value = '  +0.55%  '
$ws.Range('E26').Value = '  -0.30%  '
$ws.Range('D27').Value = '''24.51'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('E28').Value = '  +15.26%  '
$ws.Range('D29').Value = '''9.36'
$ws.Range('E29').Value = '  +2.62%  '
$ws.Range('D30').Value = '''32.15'
$ws.Range('E30').Value = '  +2.66%  '
$ws.Range('D31').Value = '''0.999'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').Value = '''17.58'
$ws.Range('E33').Value = '  +1.91%  '
$ws.Range('E34').Value = '  +8.29%  '
$ws.Range('D35').Value = '''0.0731'
$ws.Range('E35').Value = '  -1.54%  '
$ws.Range('D36').Value = '''128.77'
$ws.Range('E36').Value = '  +2.74%  '
$ws.Range('E37').Value = '  +1.13%  '
$ws.Range('E38').Value = '  +3.30%  '
$ws.Range('D39').Value = '''4.32'
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('E40').Value = '  -2.70%  '
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').Value = '''20.70'
$ws.Range('E42').Value = '  -7.39%  '
$ws.Range('D43').Value = '1.933.03'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('D44').Value = '''0.0279'
$ws.Range('E44').Value = '  -0.93%  '
$ws.Range('E45').Value = '  +2.55%  '
$ws.Range('D46').Value = '''9.23'
$ws.Range('E46').Value = '  -9.24%  '
$ws.Range('D47').Value = '''2.73'
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('D48').Value = '2.597.93'
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('E49').Value = '  +3.25%  '
$ws.Range('D50').Value = '''71.46'
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('D51').Value = '''51.42'
$ws.Range('E51').Value = '  -2.45%  '
